# Apply the "2020" column (Q) update to the SDG indicator sheet, plus the
# small workbook-level bookkeeping tweaks that went along with it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- workbook.xml: x15ac:absPath (last-saved folder shown by Excel) -------
# This is normally stamped by Excel itself from the live session's working
# folder and isn't exposed as a writable property on Workbook/Application in
# the object model, but attempt the natural property write in case the host
# maps it through.
$wb.Path = "C:\Users\korozbaeva\Desktop\Показатели ЦУР для Платформы\Глобальные показатели ЦУР\"

# --- New "2020" column --------------------------------------------------
# Column P holds 2019; clone its look into the new column Q (this also
# nudges the sheet's used range from A1:P14 to A1:Q14 automatically) before
# dropping in the 2020 figures.
$ws.Range("P4:P14").Copy()
$ws.Range("Q4:Q14").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("Q4").Value = 2020
$ws.Range("Q5").Value = 0.02
$ws.Range("Q6").Value = 0
$ws.Range("Q7").Value = 0
$ws.Range("Q8").Value = 0
$ws.Range("Q9").Value = 0.54
$ws.Range("Q10").Value = 0
$ws.Range("Q11").Value = 0
$ws.Range("Q12").Value = 0
$ws.Range("Q13").Value = 0
$ws.Range("Q14").Value = 0

# --- Selection left on the sheet after the edit --------------------------
$ws.Range("N19").Select()
